$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Int" column (column B, header "Int") is no longer produced by the
# Excel type provider's DataTypes sample - drop it entirely so Float,
# Boolean and Date shift left into B, C, D.
$ws.Range("B1").EntireColumn.Delete()

# Leave the selection where Excel naturally lands after deleting a column.
$ws.Range("B1:B1048576").Select()
